$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on D2:E51 so numeric-looking strings are not
# auto-converted to numbers (preserves literal text like "29.327.53",
# "5.225", padded percentages, etc.)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.327.53"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.861.82"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "0.7034"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Value = "238.17"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "0.07856"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("D9").Value = "0.3053"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "24.73"
$ws.Range("E10").Value = "  +6.16%  "
$ws.Range("D11").Value = "0.08163"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "1.877.61"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Value = "5.225"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").Value = "0.7133"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "89.22"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "29.396.59"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "5.820"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").Value = "0.000007788"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "238.92"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").Value = "13.20"
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").Value = "2.144.79"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "7.530"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").Value = "162.46"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").Value = "8.901"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").Value = "0.1425"
$ws.Range("E27").Value = "  -2.55%  "
$ws.Range("D28").Value = "18.09"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "1.905"
$ws.Range("E29").Value = "  -5.27%  "
$ws.Range("D30").Value = "1.373"
$ws.Range("E30").Value = "  -4.37%  "
$ws.Range("D31").Value = "1.474"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").Value = "4.301"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").Value = "4.041"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").Value = "0.05172"
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").Value = "1.181"
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("D36").Value = "0.7058"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("D37").Value = "1.003"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").Value = "2.677"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "0.01842"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "2.693"
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("D41").Value = "1.173.71"
$ws.Range("E41").Value = "  +2.76%  "
$ws.Range("D42").Value = "0.9192"
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").Value = "6.027"
$ws.Range("E43").Value = "  +1.63%  "
$ws.Range("D44").Value = "71.84"
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("D45").Value = "0.4249"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "101.83"
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").Value = "0.5352"
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("D49").Value = "1.754"
$ws.Range("E49").Value = "  -2.48%  "
$ws.Range("D50").Value = "9.156"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").Value = "6.980"
$ws.Range("E51").Value = "  -0.26%  "

# Restore the original (default) style on the range so no stray
# per-cell style index is introduced by the temporary Text format.
$ws.Range("D2:E51").Style = "Normal"
